# Applies the betexplorer "brazil serie-a 2023" update:
#  1) Re-sorts the match-detail columns (F:V) of several existing rows that
#     share the same match date but were in the wrong chronological (K-column)
#     order. Columns A:E (Indice/pais/torneio/temporada/data_partida) are left
#     untouched on every row - only F:V (home..url_partida) are rotated among
#     the rows of each affected group.
#  2) Appends one brand-new match row (row 238 / Indice 237): Corinthians 1 x 0
#     Botafogo RJ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: rotate/swap the F:V ("home" .. "url_partida") values among the rows
# listed in each cycle. For a cycle @(r0, r1, ..., rn) the rule is:
#     new(r_i) = old(r_(i+1 mod n))
# i.e. every row in the cycle receives the F:V values that used to belong to
# the next row in the cycle (wrapping around). This reproduces plain 2-row
# swaps as well as the longer 3/5-row rotations seen in the diff.
# ---------------------------------------------------------------------------

$cycles = @(
    @(63,64,65,66,67),
    @(74,75),
    @(83,85),
    @(87,88),
    @(100,101),
    @(102,103),
    @(105,106),
    @(109,110),
    @(118,119,120),
    @(132,133),
    @(134,135),
    @(155,157,156),
    @(160,162,161),
    @(163,164),
    @(182,184),
    @(210,211),
    @(216,218),
    @(219,220),
    @(221,222)
)

$firstCol = 6   # column F
$lastCol  = 22  # column V

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot the current (pre-edit) F:V values for every row in this cycle
    # before any of them get overwritten.
    $saved = @{}
    foreach ($r in $cycle) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
        }
        $saved[$r] = $rowVals
    }

    # Write each row's new content: the snapshot taken from the next row in
    # the cycle (wrapping around to the start).
    for ($i = 0; $i -lt $n; $i++) {
        $target = $cycle[$i]
        $source = $cycle[($i + 1) % $n]
        $rowVals = $saved[$source]
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($target, $c).Value = $rowVals[$c]
        }
    }
}

# ---------------------------------------------------------------------------
# Step 2: append the new match row (row 238, Indice 237).
# ---------------------------------------------------------------------------

$newRow = 238
$ws.Cells.Item($newRow, 1).Value  = 237                      # A - Indice
$ws.Cells.Item($newRow, 2).Value  = "brazil"                 # B - pais
$ws.Cells.Item($newRow, 3).Value  = "serie-a"                # C - torneio
$ws.Cells.Item($newRow, 4).Value  = "2023"                   # D - temporada
$ws.Cells.Item($newRow, 5).Value  = 45192.04166666666        # E - data_partida
$ws.Cells.Item($newRow, 6).Value  = "Corinthians"             # F - home
$ws.Cells.Item($newRow, 7).Value  = 1                         # G - home_ft_gols
$ws.Cells.Item($newRow, 8).Value  = "Botafogo RJ"             # H - away
$ws.Cells.Item($newRow, 9).Value  = 0                         # I - away_ft_gols
$ws.Cells.Item($newRow, 10).Value = 2.6                       # J - home_opening_odds
$ws.Cells.Item($newRow, 11).Value = "19/09/2023 01:12"        # K - home_opening_data_hora
$ws.Cells.Item($newRow, 12).Value = 2.91                      # L - home_closing_odds
$ws.Cells.Item($newRow, 13).Value = "23/09/2023 00:51"        # M - home_closing_data_hora
$ws.Cells.Item($newRow, 14).Value = 3.11                      # N - draw_opening_odds
$ws.Cells.Item($newRow, 15).Value = "19/09/2023 01:12"        # O - draw_opening_data_hora
$ws.Cells.Item($newRow, 16).Value = 2.99                      # P - draw_closing_odds
$ws.Cells.Item($newRow, 17).Value = "23/09/2023 00:51"        # Q - draw_closing_data_hora
$ws.Cells.Item($newRow, 18).Value = 3.01                      # R - away_opening_odds
$ws.Cells.Item($newRow, 19).Value = "19/09/2023 01:12"        # S - away_opening_data_hora
$ws.Cells.Item($newRow, 20).Value = 2.86                      # T - away_closing_odds
$ws.Cells.Item($newRow, 21).Value = "23/09/2023 00:51"        # U - away_closing_data_hora
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/brazil/serie-a/corinthians-botafogo-rj/dQOuasPr/" # V - url_partida
